$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.688.26"
$ws.Range("E2").Value = "  -0.23%  "
$ws.Range("D3").Value = "2.470.95"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "91.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.26%  "
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -1.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.90"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0850"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.42%  "
$ws.Range("E12").Value = "  -0.91%  "
$ws.Range("D13").Value = "2.851.63"
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.42"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.99%  "
$ws.Range("D16").Value = "2.474.43"
$ws.Range("E16").Value = "  -0.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.790"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.25%  "
$ws.Range("D18").Value = "41.590.08"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.43"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.22%  "
$ws.Range("D20").Value = "0.0₃0939"
$ws.Range("E20").Value = "  -1.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "239.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("E24").Value = "  +1.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.68%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("E28").Value = "  -1.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.36"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "156.85"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.43"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.04%  "
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0764"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.25%  "
$ws.Range("E35").Value = "  -0.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.95%  "
$ws.Range("E37").Value = "  +0.89%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.89"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.30%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.83"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.24%  "
$ws.Range("E40").Value = "  -0.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("E42").Value = "  -2.59%  "
$ws.Range("D43").Value = "2.001.78"
$ws.Range("E43").Value = "  +1.37%  "
$ws.Range("E44").Value = "  -0.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.66"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.81%  "
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.26%  "
$ws.Range("D48").Value = "2.731.66"
$ws.Range("E48").Value = "  +1.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "97.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "75.59"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "66.81"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.75%  "
